# Fix Training Data Issue: the "Date" column (BF) for every data row had been
# populated with "6-1-2011-12" (the source workbook's file-name-derived label)
# instead of the actual game date. Correct it to "2012-06-01" wherever that
# placeholder value is found.
#
# NOTE: assigning a date-look-alike string straight to .Value / .Value2 (or to
# .Formula as a bare literal) makes the Excel COM layer auto-convert it into a
# real date serial number (e.g. "2012-06-01" -> 41061). To keep it as literal
# text - matching the original cell's plain-string content and leaving
# NumberFormat/style untouched - we stage it as a string-returning formula
# ( ="2012-06-01" ) and then flatten just the touched cells to static values
# via Copy / PasteSpecial(xlPasteValues).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "6-1-2011-12"
$newValue = "2012-06-01"

$col = 58  # column BF
$lastRow = $ws.Cells(1, $col).End(-4121).Row  # xlDown from the header row
if ($lastRow -lt 2) { $lastRow = 31 }

$changed = $null
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Formula = '="' + $newValue + '"'
        if ($changed -eq $null) {
            $changed = $cell
        } else {
            $changed = $excel.Union($changed, $cell)
        }
    }
}

if ($changed -ne $null) {
    $changed.Copy()
    $changed.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = 0
}
